$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data block (row 164), shifting the
# existing rows 164:229 down to 168:233.
$ws.Range("A164:A167").EntireRow.Insert()

# Row 164: Crespo record / Primera
$ws.Cells.Item(164, 1).Value = 11
$ws.Cells.Item(164, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(164, 3).Value = "Bíobío"
$ws.Cells.Item(164, 4).Value = 44510
$ws.Cells.Item(164, 5).Value = 8
$ws.Cells.Item(164, 6).Value = 100112006
$ws.Cells.Item(164, 7).Value = "Repollo"
$ws.Cells.Item(164, 8).Value = "Crespo record"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 1000
$ws.Cells.Item(164, 11).Value = 700
$ws.Cells.Item(164, 12).Value = 800
$ws.Cells.Item(164, 13).Value = 750
$ws.Cells.Item(164, 14).Value = "$/unidad"
$ws.Cells.Item(164, 15).Value = "Región Metropolitana"
$ws.Cells.Item(164, 16).Value = 750
$ws.Cells.Item(164, 17).Value = 1
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# Row 165: Crespo record / Segunda
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = 44510
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = 100112006
$ws.Cells.Item(165, 7).Value = "Repollo"
$ws.Cells.Item(165, 8).Value = "Crespo record"
$ws.Cells.Item(165, 9).Value = "Segunda"
$ws.Cells.Item(165, 10).Value = 500
$ws.Cells.Item(165, 11).Value = 600
$ws.Cells.Item(165, 12).Value = 600
$ws.Cells.Item(165, 13).Value = 600
$ws.Cells.Item(165, 14).Value = "$/unidad"
$ws.Cells.Item(165, 15).Value = "Región Metropolitana"
$ws.Cells.Item(165, 16).Value = 600
$ws.Cells.Item(165, 17).Value = 1
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# Row 166: Morada(o) / Primera
$ws.Cells.Item(166, 1).Value = 11
$ws.Cells.Item(166, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(166, 3).Value = "Bíobío"
$ws.Cells.Item(166, 4).Value = 44510
$ws.Cells.Item(166, 5).Value = 8
$ws.Cells.Item(166, 6).Value = 100112006
$ws.Cells.Item(166, 7).Value = "Repollo"
$ws.Cells.Item(166, 8).Value = "Morada(o)"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 600
$ws.Cells.Item(166, 11).Value = 800
$ws.Cells.Item(166, 12).Value = 900
$ws.Cells.Item(166, 13).Value = 850
$ws.Cells.Item(166, 14).Value = "$/unidad"
$ws.Cells.Item(166, 15).Value = "Región Metropolitana"
$ws.Cells.Item(166, 16).Value = 850
$ws.Cells.Item(166, 17).Value = 1
$ws.Cells.Item(166, 18).Value = "Hortaliza"

# Row 167: Morada(o) / Segunda
$ws.Cells.Item(167, 1).Value = 11
$ws.Cells.Item(167, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(167, 3).Value = "Bíobío"
$ws.Cells.Item(167, 4).Value = 44510
$ws.Cells.Item(167, 5).Value = 8
$ws.Cells.Item(167, 6).Value = 100112006
$ws.Cells.Item(167, 7).Value = "Repollo"
$ws.Cells.Item(167, 8).Value = "Morada(o)"
$ws.Cells.Item(167, 9).Value = "Segunda"
$ws.Cells.Item(167, 10).Value = 300
$ws.Cells.Item(167, 11).Value = 700
$ws.Cells.Item(167, 12).Value = 700
$ws.Cells.Item(167, 13).Value = 700
$ws.Cells.Item(167, 14).Value = "$/unidad"
$ws.Cells.Item(167, 15).Value = "Región Metropolitana"
$ws.Cells.Item(167, 16).Value = 700
$ws.Cells.Item(167, 17).Value = 1
$ws.Cells.Item(167, 18).Value = "Hortaliza"
